$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: set new values
$ws.Range("C23:F23").Value = 5

# Apply style (no-fill bordered style) matching G23..K23 by copying format
$src = $ws.Range("G23")
$dst = $ws.Range("C23:F23")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Clear G23 entirely (value + style) -> cell disappears from the XML
$ws.Range("G23").Clear()

# Formula for L23 explicit (breaks shared group for this cell only)
$ws.Range("L23").Formula = "=SUM(C23:F23)"

Write-Host "done"
